$wb = $excel.ActiveWorkbook

# Duplicate the "Czech_MXPanel" sheet (last sheet) to create the new "Swiss" sheet,
# then update its content for the Switzerland market.
$source = $wb.Worksheets.Item("Czech_MXPanel")
$source.Copy($null, $source)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Swiss"

$ws.Range("B2").Value = "Switzerland Market"
$ws.Range("B4").Value = "NGC-3476/T2344"
$ws.Range("A11").Value = "PROFILE Communicator"

# Update leftover selection artifacts on the other sheets, mirroring the
# navigation that happened while adding the new market data.
$czech = $wb.Worksheets.Item("Czech")
$czech.Activate() | Out-Null
$czech.Range("A7:A17").Select() | Out-Null

$czechMx = $wb.Worksheets.Item("Czech_MXPanel")
$czechMx.Activate() | Out-Null
$czechMx.Cells.Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("B9").Select() | Out-Null
